$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.424.85"
Set-TextValue $ws.Range("E2") "  -0.94%  "
Set-TextValue $ws.Range("D3") "3.528.08"
Set-TextValue $ws.Range("E3") "  +0.43%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.22%  "
Set-TextValue $ws.Range("D5") "579.00"
Set-TextValue $ws.Range("E5") "  +5.56%  "
Set-TextValue $ws.Range("D6") "179.28"
Set-TextValue $ws.Range("E6") "  -5.78%  "
Set-TextValue $ws.Range("E7") "  +4.39%  "
Set-TextValue $ws.Range("E8") "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.639"
Set-TextValue $ws.Range("E9") "  +1.30%  "
Set-TextValue $ws.Range("D10") "0.162"
Set-TextValue $ws.Range("E10") "  +7.76%  "
Set-TextValue $ws.Range("D11") "55.73"
Set-TextValue $ws.Range("E11") "  +1.89%  "
Set-TextValue $ws.Range("D12") "0.0000280"
Set-TextValue $ws.Range("E12") "  +4.59%  "
Set-TextValue $ws.Range("D13") "9.30"
Set-TextValue $ws.Range("E13") "  -0.40%  "
Set-TextValue $ws.Range("D14") "4.087.68"
Set-TextValue $ws.Range("E14") "  +0.28%  "
Set-TextValue $ws.Range("D15") "3.531.22"
Set-TextValue $ws.Range("E15") "  +0.46%  "
Set-TextValue $ws.Range("E16") "  +0.32%  "
Set-TextValue $ws.Range("D17") "18.45"
Set-TextValue $ws.Range("E17") "  +1.81%  "
Set-TextValue $ws.Range("D18") "66.335.35"
Set-TextValue $ws.Range("E18") "  -1.16%  "
Set-TextValue $ws.Range("E19") "  +1.63%  "
Set-TextValue $ws.Range("E20") "  +1.75%  "
Set-TextValue $ws.Range("D21") "416.49"
Set-TextValue $ws.Range("E21") "  -1.92%  "
Set-TextValue $ws.Range("D22") "4.23"
Set-TextValue $ws.Range("E22") "  +8.28%  "
Set-TextValue $ws.Range("D23") "4.32"
Set-TextValue $ws.Range("E23") "  +3.54%  "
Set-TextValue $ws.Range("D24") "86.02"
Set-TextValue $ws.Range("E24") "  +1.75%  "
Set-TextValue $ws.Range("D25") "13.20"
Set-TextValue $ws.Range("E25") "  +10.73%  "
Set-TextValue $ws.Range("D26") "11.35"
Set-TextValue $ws.Range("E26") "  +1.96%  "
Set-TextValue $ws.Range("E27") "  -0.79%  "
Set-TextValue $ws.Range("D28") "6.05"
Set-TextValue $ws.Range("E28") "  -1.74%  "
Set-TextValue $ws.Range("D29") "9.15"
Set-TextValue $ws.Range("E29") "  +3.55%  "
Set-TextValue $ws.Range("D30") "30.59"
Set-TextValue $ws.Range("E30") "  +1.65%  "
Set-TextValue $ws.Range("D33") "11.74"
Set-TextValue $ws.Range("E33") "  +0.52%  "
Set-TextValue $ws.Range("E34") "  +1.16%  "
Set-TextValue $ws.Range("E35") "  +12.42%  "
Set-TextValue $ws.Range("D36") "59.86"
Set-TextValue $ws.Range("E36") "  +0.74%  "
Set-TextValue $ws.Range("D37") "0.0₃0812"
Set-TextValue $ws.Range("E37") "  -0.01%  "
Set-TextValue $ws.Range("D38") "1.00"
Set-TextValue $ws.Range("E38") "  +0.17%  "
Set-TextValue $ws.Range("D39") "37.33"
Set-TextValue $ws.Range("E39") "  -2.84%  "
Set-TextValue $ws.Range("D40") "3.55"
Set-TextValue $ws.Range("E40") "  +7.42%  "
Set-TextValue $ws.Range("E41") "  -0.60%  "
Set-TextValue $ws.Range("D42") "3.251.86"
Set-TextValue $ws.Range("E42") "  +8.23%  "
Set-TextValue $ws.Range("D43") "1.00"
Set-TextValue $ws.Range("E43") "  +0.09%  "
Set-TextValue $ws.Range("D44") "2.94"
Set-TextValue $ws.Range("E44") "  +2.89%  "
Set-TextValue $ws.Range("E45") "  -3.04%  "
Set-TextValue $ws.Range("D46") "0.0422"
Set-TextValue $ws.Range("E46") "  +1.65%  "
Set-TextValue $ws.Range("E47") "  -2.14%  "
Set-TextValue $ws.Range("E48") "  -0.33%  "
Set-TextValue $ws.Range("E49") "  +2.24%  "
Set-TextValue $ws.Range("D50") "8.66"
Set-TextValue $ws.Range("E50") "  -0.56%  "
Set-TextValue $ws.Range("D51") "138.55"
Set-TextValue $ws.Range("E51") "  -1.12%  "

# Row 31/32 swap: Bittensor <-> NEARProtocol (ranked by position; full row content swap + updated values)
Set-TextValue $ws.Range("B31") "NEARProtocol"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D31") "6.63"
Set-TextValue $ws.Range("E31") "  -0.11%  "

Set-TextValue $ws.Range("B32") "Bittensor"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D32") "617.34"
Set-TextValue $ws.Range("E32") "  -5.44%  "
